# Menu-Languages.docx (Estonian) - "New translations" commit.
#
# The only meaningful content change is the localized label for the
# "RPC Explorer" menu entry (word/document.xml), which was renamed to
# "Insight Explorer".
#
# A plain Find/Replace performs the textual substitution correctly, but
# this engine rebuilds the <w:t> element without xml:space="preserve" once
# the replacement text has no leading/trailing whitespace to justify it -
# the canonical (target) OOXML keeps that attribute on the run (it is
# simply inherited unchanged from the original run). To reproduce the
# target byte-for-byte we rebuild the containing paragraph explicitly via
# Range.InsertXML, keeping xml:space="preserve" on the new text.

$d = $word.ActiveDocument

$oldText = "RPC Explorer"
$newText = "Insight Explorer"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = $oldText
$find.Forward = $true
$find.Wrap = 1

$replacedPrecisely = $false

if ($find.Execute()) {
    $hit = $find.Parent.Duplicate
    $para = $hit.Paragraphs(1)
    $paraRange = $para.Range

    # Only take the exact-XML-rebuild path when the paragraph still has the
    # exact shape we captured from the source document (a tab run followed
    # by the "RPC Explorer" text run, optionally followed by Word's
    # trailing paragraph-mark character) - otherwise fall back below.
    $paraPlainText = $paraRange.Text.TrimEnd([char]13, [char]7)
    if ($paraPlainText -eq "`t$oldText") {
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' +
            '<w:p>' +
            '<w:pPr>' +
            '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
            '<w:ind w:left="360"/>' +
            '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:color w:val="000000"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '</w:rPr>' +
            '</w:pPr>' +
            '<w:r>' +
            '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>' +
            '<w:color w:val="4472C4"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
            '</w:rPr>' +
            '<w:tab/>' +
            '</w:r>' +
            '<w:r>' +
            '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:eastAsia="Times New Roman"/>' +
            '<w:color w:val="000000"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '</w:rPr>' +
            "<w:t xml:space=`"preserve`">$newText</w:t>" +
            '</w:r>' +
            '</w:p>' +
            '</w:body>' +
            '</w:document>' +
            '</pkg:xmlData>' +
            '</pkg:part>' +
            '</pkg:package>'

        $paraRange.InsertXML($xml) | Out-Null
        $replacedPrecisely = $true
    }
}

if (-not $replacedPrecisely) {
    # Fallback: plain text replacement (still correct content, just without
    # the redundant xml:space="preserve" nuance on the rebuilt run).
    $find2 = $d.Content.Find
    $find2.ClearFormatting()
    $find2.Replacement.ClearFormatting()
    $find2.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}
